$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.243.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.269.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.84%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "113.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "264.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.620"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -3.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.05%  "
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("E14").Value = "  -2.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.607.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.856"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.266.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.112.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.63%  "
$ws.Range("E19").Value = "  -3.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "230.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.42%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.46%  "
$ws.Range("E29").Value = "  -1.23%  "
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("E31").Value = "  -1.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "171.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.47%  "
$ws.Range("E33").Value = "  -3.14%  "
$ws.Range("E34").Value = "  -1.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.64"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("E37").Value = "  -4.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0350"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.50%  "
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("E40").Value = "  -8.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +15.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.21%  "
$ws.Range("E43").Value = "  +3.36%  "
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.91%  "
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("E47").Value = "  -2.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.33%  "
$ws.Range("E49").Value = "  -3.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "100.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.15%  "
$ws.Range("E51").Value = "  +0.76%  "
